# "after logout test done"
#
# Login Test sheet: the existing admin/prison login case now records a
# negative ("No") result, and four new login attempts (all against the
# "admin_welikada" account, exercising a correct password plus three bad
# variants - leading space, trailing space, and wrong case) are appended.
# Selection/active-tab view state moves to the Login Test sheet, and the
# previously-active Allocate Location Test sheet's selection is nudged to
# B3. A hidden helper defined name (LOCAL_MYSQL_DATE_FORMAT) is also added
# to the workbook.

$wb = $excel.ActiveWorkbook

# --- Login Test sheet: new rows + updated result ------------------------
$wsLogin = $wb.Worksheets.Item("Login Test")

$wsLogin.Range("A3").Value = "admin_welikada"
$wsLogin.Range("B3").Value = "test`$123"

$wsLogin.Range("A4").Value = "admin_welikada"
$wsLogin.Range("B4").Value = " test`$123"

$wsLogin.Range("A5").Value = "admin_welikada"
$wsLogin.Range("B5").Value = "test`$123 "

$wsLogin.Range("A6").Value = "admin_welikada"
$wsLogin.Range("B6").Value = "TEST`$123"

$wsLogin.Range("C3").Value = "Yes"
$wsLogin.Range("C4").Value = "No"
$wsLogin.Range("C5").Value = "No"
$wsLogin.Range("C6").Value = "No"

# Original admin/prison run is now recorded as a failed ("No") test case.
$wsLogin.Range("C2").Value = "No"

# --- Workbook-level hidden defined name ----------------------------------
$mysqlDateFormat = '=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&" "&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)'
$dateFormatName = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", $mysqlDateFormat)
$dateFormatName.Visible = $false

# --- View state: Allocate Location Test selection moves to B3 -----------
$wsAllocate = $wb.Worksheets.Item("Allocate Location Test")
$wsAllocate.Range("B3").Select()

# --- View state: Login Test becomes the active tab, selection at E5 -----
$wsLogin.Activate()
$wsLogin.Range("E5").Select()
